$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the record, keyed by letter for readability.
# (Y/AA - the start/end dates - are intentionally excluded: every
# affected row shares the same date, so swapping them is a no-op that
# would only risk an unwanted type change on those cells.)
$cols = @{
    A = 1; B = 2; D = 4; E = 5; F = 6; G = 7; H = 8; I = 9; P = 16;
    Q = 17; R = 18; S = 19; T = 20; U = 21; V = 22; W = 23;
    Z = 26; AB = 28; AD = 30; AE = 31; AG = 33; AT = 46;
    AW = 49; AX = 50; AY = 51
}

function Swap-Rows($rowA, $rowB) {
    foreach ($col in $cols.Values) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valueA = $cellA.Value2
        $valueB = $cellB.Value2
        $cellA.Value2 = $valueB
        $cellB.Value2 = $valueA
    }
}

# The three record pairs that were reordered in the source data.
Swap-Rows 8 9
Swap-Rows 11 12
Swap-Rows 13 14
